$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.908.02"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "2.626.97"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "2.625.16"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("E10").Value = "  +14.44%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000190"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.31%  "
$ws.Range("D16").Value = "3.096.44"
$ws.Range("D17").Value = "67.788.01"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").Value = "2.626.76"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +3.82%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "585.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +5.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").Value = "0.0₆0292"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.626"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.55%  "
